$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Series" row describing MCH204-1 (identifier, level, extent, notes)
$ws.Range("A2").Value = "MCH204-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24C | GRAP COUNT NUMER: NONE"

# Match the font used for the new row's cells (Calibri 10pt, automatic/theme text colour).
# Touching C2/D2/H2 here (even with no value) materialises them as blank styled cells.
foreach ($addr in @("A2", "C2", "D2", "E2", "F2", "G2", "H2")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.Font.ThemeColor = 1
}

$ws.Range("A2:H2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
